{"js": "// Update the worksheet date and the 25 two-digit multiplication problems\n// to the new set of values. Every string below is unique within the\n// document, so a plain search-and-replace is unambiguous.\nconst replacements = [\n  [\"2023-10-14 Saturday\", \"2023-10-15 Sunday\"],\n  [\"23\u00d736=828\", \"97\u00d753=5141\"],\n  [\"73\u00d791=6643\", \"25\u00d799=2475\"],\n  [\"37\u00d712=444\", \"49\u00d799=4851\"],\n  [\"79\u00d711=869\", \"46\u00d732=1472\"],\n  [\"87\u00d727=2349\", \"95\u00d734=3230\"],\n  [\"58\u00d789=5162\", \"51\u00d788=4488\"],\n  [\"80\u00d733=2640\", \"40\u00d748=1920\"],\n  [\"91\u00d711=1001\", \"75\u00d778=5850\"],\n  [\"96\u00d737=3552\", \"77\u00d718=1386\"],\n  [\"69\u00d746=3174\", \"92\u00d779=7268\"],\n  [\"59\u00d773=4307\", \"30\u00d757=1710\"],\n  [\"22\u00d754=1188\", \"60\u00d772=4320\"],\n  [\"12\u00d796=1152\", \"88\u00d718=1584\"],\n  [\"13\u00d768=884\", \"43\u00d736=1548\"],\n  [\"62\u00d732=1984\", \"24\u00d736=864\"],\n  [\"11\u00d742=462\", \"65\u00d719=1235\"],\n  [\"71\u00d772=5112\", \"84\u00d786=7224\"],\n  [\"23\u00d767=1541\", \"30\u00d714=420\"],\n  [\"73\u00d761=4453\", \"63\u00d728=1764\"],\n  [\"21\u00d726=546\", \"96\u00d734=3264\"],\n  [\"39\u00d721=819\", \"74\u00d750=3700\"],\n  [\"52\u00d751=2652\", \"56\u00d763=3528\"],\n  [\"67\u00d756=3752\", \"19\u00d737=703\"],\n  [\"81\u00d736=2916\", \"56\u00d782=4592\"],\n  [\"43\u00d715=645\", \"81\u00d776=6156\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "# Update the worksheet date and the 25 two-digit multiplication problems\n# to the new set of values. Every string below is unique within the\n# document, so a plain Find/Replace is unambiguous.\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @(\"2023-10-14 Saturday\", \"2023-10-15 Sunday\"),\n    @(\"23\u00d736=828\", \"97\u00d753=5141\"),\n    @(\"73\u00d791=6643\", \"25\u00d799=2475\"),\n    @(\"37\u00d712=444\", \"49\u00d799=4851\"),\n    @(\"79\u00d711=869\", \"46\u00d732=1472\"),\n    @(\"87\u00d727=2349\", \"95\u00d734=3230\"),\n    @(\"58\u00d789=5162\", \"51\u00d788=4488\"),\n    @(\"80\u00d733=2640\", \"40\u00d748=1920\"),\n    @(\"91\u00d711=1001\", \"75\u00d778=5850\"),\n    @(\"96\u00d737=3552\", \"77\u00d718=1386\"),\n    @(\"69\u00d746=3174\", \"92\u00d779=7268\"),\n    @(\"59\u00d773=4307\", \"30\u00d757=1710\"),\n    @(\"22\u00d754=1188\", \"60\u00d772=4320\"),\n    @(\"12\u00d796=1152\", \"88\u00d718=1584\"),\n    @(\"13\u00d768=884\", \"43\u00d736=1548\"),\n    @(\"62\u00d732=1984\", \"24\u00d736=864\"),\n    @(\"11\u00d742=462\", \"65\u00d719=1235\"),\n    @(\"71\u00d772=5112\", \"84\u00d786=7224\"),\n    @(\"23\u00d767=1541\", \"30\u00d714=420\"),\n    @(\"73\u00d761=4453\", \"63\u00d728=1764\"),\n    @(\"21\u00d726=546\", \"96\u00d734=3264\"),\n    @(\"39\u00d721=819\", \"74\u00d750=3700\"),\n    @(\"52\u00d751=2652\", \"56\u00d763=3528\"),\n    @(\"67\u00d756=3752\", \"19\u00d737=703\"),\n    @(\"81\u00d736=2916\", \"56\u00d782=4592\"),\n    @(\"43\u00d715=645\", \"81\u00d776=6156\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n    if (-not $ok) {\n        throw \"Text not found: $oldText\"\n    }\n}"}
